$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3249")

# Update the activation date (shared by "Ativação:" row and the mis-linked
# cells in the "Programa:" row) from 01/01/2012 to 01/01/2023.
# Force text storage (not an auto-converted date serial number).
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"

$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("B15").Value = "01/01/2023"
$ws.Range("C15").Value = "01/01/2023"

# Fill in the previously empty "Objectives:" detail cells.
$ws.Range("B11").Value = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."
$ws.Range("C11").Value = "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."

# Fill in the previously empty "Short syllabus:" detail cells.
$ws.Range("B14").Value = "To be defined, according to the programmed topic."
$ws.Range("C14").Value = "To be defined, according to the programmed topic."

# Fill in the previously empty "Syllabus:" detail cells.
$ws.Range("B16").Value = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C16").Value = "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."

# Match the existing column styles used throughout the sheet: column B uses
# the normal wrap-text style, column C uses the red "modified" wrap-text style.
$ws.Range("B11").Style = $ws.Range("B13").Style
$ws.Range("C11").Style = $ws.Range("C13").Style

$ws.Range("B14").Style = $ws.Range("B13").Style
$ws.Range("C14").Style = $ws.Range("C13").Style

$ws.Range("B16").Style = $ws.Range("B13").Style
$ws.Range("C16").Style = $ws.Range("C13").Style
